# New weekly price observation for "Coco" was added to the Vega Modelo de
# Temuco daily-logic subset. Insert a brand-new row at position 14 (pushing
# the existing rows 14-97 down to 15-98, and extending the used range from
# A1:T97 to A1:T98), then populate the new row with the reported values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("14:14").Insert()

$ws.Range("A14").Value = 10
$ws.Range("B14").Value = "Vega Modelo de Temuco"
$ws.Range("C14").Value = "La Araucanía"
$ws.Range("D14").Value = 44901
$ws.Range("E14").Value = 9
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100108
$ws.Range("H14").Value = "Tropicales y subtropicales"
$ws.Range("I14").Value = 100108007
$ws.Range("J14").Value = "Coco"
$ws.Range("K14").Value = "Sin especificar"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 20
$ws.Range("N14").Value = 30000
$ws.Range("O14").Value = 30000
$ws.Range("P14").Value = 30000
$ws.Range("Q14").Value = "$/malla 20 unidades"
$ws.Range("R14").Value = "Perú"
$ws.Range("S14").Value = 1500
$ws.Range("T14").Value = 20
